$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.897.27'
$ws.Range("E2").Value = '  +4.81%  '

$ws.Range("D3").Value = '3.527.65'
$ws.Range("E3").Value = '  +7.28%  '

$ws.Range("E4").Value = '  +0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '188.96'
$ws.Range("E5").Value = '  +8.14%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '559.75'
$ws.Range("E6").Value = '  +4.67%  '

$ws.Range("D7").Value = '3.522.76'
$ws.Range("E7").Value = '  +7.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.615'
$ws.Range("E8").Value = '  +2.69%  '

$ws.Range("E9").Value = '  +0.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.633'
$ws.Range("E10").Value = '  +3.21%  '

$ws.Range("E11").Value = '  +11.41%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.77'
$ws.Range("E12").Value = '  +1.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000270'
$ws.Range("E13").Value = '  +3.76%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.40'
$ws.Range("E14").Value = '  +1.94%  '

$ws.Range("D15").Value = '4.100.97'
$ws.Range("E15").Value = '  +8.60%  '

$ws.Range("D16").Value = '3.535.97'
$ws.Range("E16").Value = '  +8.74%  '

$ws.Range("E17").Value = '  +3.46%  '

$ws.Range("D18").Value = '66.840.32'
$ws.Range("E18").Value = '  +5.26%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.23'
$ws.Range("E19").Value = '  +4.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.02'
$ws.Range("E20").Value = '  +7.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.995'
$ws.Range("E21").Value = '  +2.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '432.59'
$ws.Range("E22").Value = '  +16.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.11'
$ws.Range("E23").Value = '  +8.48%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.29'
$ws.Range("E24").Value = '  +4.81%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.13'
$ws.Range("E25").Value = '  +0.98%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.07'
$ws.Range("E26").Value = '  -2.45%  '

$ws.Range("E27").Value = '  +8.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.23'
$ws.Range("E28").Value = '  +7.20%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.13'
$ws.Range("E29").Value = '  +9.79%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.46'
$ws.Range("E30").Value = '  +5.65%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '644.74'
$ws.Range("E31").Value = '  +0.30%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.57'
$ws.Range("E32").Value = '  +0.32%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.73'
$ws.Range("E33").Value = '  +3.41%  '

$ws.Range("E34").Value = '  +3.38%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '59.61'
$ws.Range("E35").Value = '  +4.35%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '38.42'
$ws.Range("E36").Value = '  +3.77%  '

$ws.Range("D37").Value = '0.0₃0808'
$ws.Range("E37").Value = '  +7.90%  '

$ws.Range("E38").Value = '  +17.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.390'
$ws.Range("E40").Value = '  +1.87%  '

$ws.Range("E41").Value = '  +12.98%  '

$ws.Range("E42").Value = '  +0.75%  '

$ws.Range("D43").Value = '3.039.44'
$ws.Range("E43").Value = '  +3.87%  '

$ws.Range("E44").Value = '  +2.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.87'
$ws.Range("E45").Value = '  +9.35%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.34'
$ws.Range("E46").Value = '  +7.95%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0418'
$ws.Range("E47").Value = '  +4.34%  '

$ws.Range("E48").Value = '  +2.69%  '

$ws.Range("E49").Value = '  +5.22%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.79'
$ws.Range("E50").Value = '  +4.82%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.62'
$ws.Range("E51").Value = '  +8.95%  '
